# Scheduled-runner refresh of market-price derived columns (H-N) across
# the per-job-class leve-profit sheets. Values mirror a re-pull of
# currentAveragePrice(NQ/HQ) data and the resulting LevePrice/LeveProfit
# recalculations for the affected rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 936.25
$ws.Range("I32").Value = 897
$ws.Range("K32").Value = 897
$ws.Range("M32").Value = -571

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H132").Value = 1169.625
$ws.Range("I132").Value = 1175.0968
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3525.2904
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -995.2903999999999
$ws.Range("N132").Value = -8060

$ws.Range("H137").Value = 1670.125
$ws.Range("I137").Value = 1297.75
$ws.Range("K137").Value = 3893.25
$ws.Range("M137").Value = -1343.25

$ws.Range("H138").Value = 2500.5312
$ws.Range("I138").Value = 2698.1904
$ws.Range("K138").Value = 8094.5712
$ws.Range("M138").Value = -2954.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3261.879
$ws.Range("I32").Value = 2109.283
$ws.Range("K32").Value = 2109.283
$ws.Range("M32").Value = -1822.283

$ws.Range("H102").Value = 1033.5
$ws.Range("I102").Value = 500
$ws.Range("J102").Value = 1567
$ws.Range("K102").Value = 500
$ws.Range("L102").Value = 1567
$ws.Range("M102").Value = 1122
$ws.Range("N102").Value = -4811

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1842.5714
$ws.Range("I94").Value = 1858
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 1858
$ws.Range("L94").Value = 1750
$ws.Range("M94").Value = -1407
$ws.Range("N94").Value = -2652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 820.86664
$ws.Range("I134").Value = 829.5714
$ws.Range("J134").Value = 699
$ws.Range("K134").Value = 2488.7142
$ws.Range("L134").Value = 2097
$ws.Range("M134").Value = 46.28579999999965
$ws.Range("N134").Value = -7167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 578.2143
$ws.Range("J5").Value = 952.5
$ws.Range("L5").Value = 2857.5
$ws.Range("N5").Value = -3081.5

$ws.Range("H7").Value = 232.3
$ws.Range("I7").Value = 77.875
$ws.Range("K7").Value = 233.625
$ws.Range("M7").Value = -121.625

$ws.Range("H17").Value = 1900
$ws.Range("I17").Value = 1900
$ws.Range("K17").Value = 5700
$ws.Range("M17").Value = -5531

$ws.Range("H34").Value = 1380
$ws.Range("J34").Value = 2000
$ws.Range("L34").Value = 6000
$ws.Range("N34").Value = -6168

$ws.Range("H39").Value = 2599.6667
$ws.Range("J39").Value = 2599.6667
$ws.Range("L39").Value = 7799.000100000001
$ws.Range("N39").Value = -8387.000100000001

$ws.Range("H68").Value = 1000
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4622

$ws.Range("H71").Value = 1000
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 9000
$ws.Range("N71").Value = -17112

$ws.Range("H92").Value = 366.66666
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H121").Value = 652.55554
$ws.Range("J121").Value = 706.6667
$ws.Range("L121").Value = 2120.0001
$ws.Range("N121").Value = -4740.0001

$ws.Range("H125").Value = 1707.5
$ws.Range("I125").Value = 1707.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 5122.5
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -202.5

$ws.Range("H129").Value = 81160.22
$ws.Range("I129").Value = 595.3333
$ws.Range("J129").Value = 121442.664
$ws.Range("K129").Value = 1785.9999
$ws.Range("L129").Value = 364327.992
$ws.Range("M129").Value = 3214.0001
$ws.Range("N129").Value = -374327.992

$ws.Range("H131").Value = 780.4299999999999
$ws.Range("J131").Value = 787.07294
$ws.Range("L131").Value = 2361.21882
$ws.Range("N131").Value = -12441.21882

$ws.Range("H135").Value = 578.2143
$ws.Range("J135").Value = 952.5
$ws.Range("L135").Value = 8572.5
$ws.Range("N135").Value = -13642.5

$ws.Range("H139").Value = 9368.923000000001
$ws.Range("I139").Value = 10709.091
$ws.Range("J139").Value = 1998
$ws.Range("K139").Value = 32127.273
$ws.Range("L139").Value = 5994
$ws.Range("M139").Value = -26987.273
$ws.Range("N139").Value = -16274

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 24999.5
$ws.Range("J49").Value = 24999.5
$ws.Range("L49").Value = 24999.5
$ws.Range("N49").Value = -25367.5

$ws.Range("H122").Value = 1712.5151
$ws.Range("I122").Value = 1530.5555
$ws.Range("J122").Value = 2531.3333
$ws.Range("K122").Value = 4591.666499999999
$ws.Range("L122").Value = 7593.999899999999
$ws.Range("M122").Value = -2141.666499999999
$ws.Range("N122").Value = -12493.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5725.7334
$ws.Range("I7").Value = 2856.2856
$ws.Range("K7").Value = 2856.2856
$ws.Range("M7").Value = -2744.2856

$ws.Range("H22").Value = 1860
$ws.Range("I22").Value = 1724
$ws.Range("K22").Value = 1724
$ws.Range("M22").Value = -1429

$ws.Range("H27").Value = 1860
$ws.Range("I27").Value = 1724
$ws.Range("K27").Value = 1724
$ws.Range("M27").Value = -1617

$ws.Range("H46").Value = 1399.6666
$ws.Range("I46").Value = 1300
$ws.Range("J46").Value = 1412.125
$ws.Range("K46").Value = 1300
$ws.Range("L46").Value = 1412.125
$ws.Range("M46").Value = -1112
$ws.Range("N46").Value = -1788.125

$ws.Range("H61").Value = 3714.2856
$ws.Range("I61").Value = 3200
$ws.Range("K61").Value = 3200
$ws.Range("M61").Value = -2998

$ws.Range("H113").Value = 3714.2856
$ws.Range("I113").Value = 3200
$ws.Range("K113").Value = 3200
$ws.Range("M113").Value = -1030

$ws.Range("H122").Value = 7830.55
$ws.Range("I122").Value = 6345.1113
$ws.Range("K122").Value = 19035.3339
$ws.Range("M122").Value = -16585.3339

$ws.Range("H126").Value = 5725.7334
$ws.Range("I126").Value = 2856.2856
$ws.Range("K126").Value = 8568.856800000001
$ws.Range("M126").Value = -6098.856800000001

$ws.Range("H137").Value = 28459
$ws.Range("J137").Value = 28459
$ws.Range("L137").Value = 28459
$ws.Range("N137").Value = -38659

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 28449.5
$ws.Range("J98").Value = 28449.5
$ws.Range("L98").Value = 28449.5
$ws.Range("N98").Value = -34439.5

$ws.Range("H126").Value = 5124.357
$ws.Range("I126").Value = 4294.25
$ws.Range("K126").Value = 12882.75
$ws.Range("M126").Value = -10412.75
